# Insert a new data row at row 7 (pushing the existing rows 7..33 down to 8..34),
# then populate the new row with this week's price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert()

$ws.Cells.Item(7, 1).Value = 11
$ws.Cells.Item(7, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value = "Bíobío"
$ws.Cells.Item(7, 4).Value = 44624
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = 100112030
$ws.Cells.Item(7, 7).Value = "Poroto granado"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 150
$ws.Cells.Item(7, 11).Value = 25000
$ws.Cells.Item(7, 12).Value = 26000
$ws.Cells.Item(7, 13).Value = 25467
$ws.Cells.Item(7, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(7, 15).Value = "Región Metropolitana"
$ws.Cells.Item(7, 16).Value = 1019
$ws.Cells.Item(7, 17).Value = 25
$ws.Cells.Item(7, 18).Value = "Hortaliza"
